$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 74, pushing the existing data (old rows 74-123)
# down to rows 76-125. Excel copies formatting (incl. the date style on column D)
# from the row above automatically.
$ws.Rows.Item(74).Insert()
$ws.Rows.Item(74).Insert()

# New row 74: same record shape as the (now shifted) row 76, but for the
# "Fortuna" variety, sold as "$/bandeja" instead of "$/caja", dated 45001.
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = "Vega Monumental Concepción"
$ws.Range("C74").Value = "Bíobío"
$ws.Range("D74").Value = 45001
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100103
$ws.Range("H74").Value = "Frutos de hueso (carozo)"
$ws.Range("I74").Value = 100103002
$ws.Range("J74").Value = "Ciruela"
$ws.Range("K74").Value = "Fortuna"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 200
$ws.Range("N74").Value = 9000
$ws.Range("O74").Value = 10000
$ws.Range("P74").Value = 9500
$ws.Range("Q74").Value = "$/bandeja 18 kilos granel"
$ws.Range("R74").Value = "Región de O'Higgins"
$ws.Range("S74").Value = 528
$ws.Range("T74").Value = 18

# New row 75: same pattern, "Segunda" quality.
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = "Vega Monumental Concepción"
$ws.Range("C75").Value = "Bíobío"
$ws.Range("D75").Value = 45001
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100103
$ws.Range("H75").Value = "Frutos de hueso (carozo)"
$ws.Range("I75").Value = 100103002
$ws.Range("J75").Value = "Ciruela"
$ws.Range("K75").Value = "Fortuna"
$ws.Range("L75").Value = "Segunda"
$ws.Range("M75").Value = 100
$ws.Range("N75").Value = 8000
$ws.Range("O75").Value = 8000
$ws.Range("P75").Value = 8000
$ws.Range("Q75").Value = "$/bandeja 18 kilos granel"
$ws.Range("R75").Value = "Región de O'Higgins"
$ws.Range("S75").Value = 444
$ws.Range("T75").Value = 18
